$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 297, pushing the existing rows 297-300 down to 301-304
$ws.Range("A297:T300").Insert()

# Populate the 4 newly inserted rows (297-300) with the new weekly data
$newRows = @(
    @{ Row = 297; D = 44628; L = "Especial"; M = 280; N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos granel"; R = "Provincia de Chacabuco"; S = 667; T = 18 },
    @{ Row = 298; D = 44628; L = "Primera";  M = 310; N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos granel"; R = "Provincia de Chacabuco"; S = 556; T = 18 },
    @{ Row = 299; D = 44628; L = "Segunda";  M = 350; N = 6000;  O = 6000;  P = 6000;  Q = "`$/caja 18 kilos granel"; R = "Provincia de Chacabuco"; S = 333; T = 18 },
    @{ Row = 300; D = 44628; L = "Tercera";  M = 200; N = 3000;  O = 3000;  P = 3000;  Q = "`$/caja 18 kilos granel"; R = "Provincia de Chacabuco"; S = 167; T = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Range("D$row").NumberFormat = $ws.Range("D296").NumberFormat
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107011
    $ws.Cells.Item($row, 10).Value = "Tuna"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
